$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: multi-line ("was X / now Y -Z%") price-drop notation
$row13 = @(
    "'`$59.90",
    "'`$169.00",
    "'`$39.90",
    "'`$45.00",
    "'`$219.00",
    "'`$24.90`n`$29.90-17%",
    "'`$851.00`n`$889.00-4%",
    "'`$851.00`n`$889.00-4%",
    "'`$851.00`n`$889.00-4%",
    "'`$289.00`n`$369.00-22%"
)

# Row 14: same columns, but plain current price only (no drop annotation)
$row14 = @(
    "'`$59.90",
    "'`$169.00",
    "'`$39.90",
    "'`$45.00",
    "'`$219.00",
    "'`$24.90",
    "'`$851.00",
    "'`$851.00",
    "'`$851.00",
    "'`$289.00"
)

# Column A holds the date serial, formatted like A12 (yyyy-mm-dd).
$ws.Cells.Item(13, 1).Value = 44207
$ws.Cells.Item(13, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(14, 1).Value = 44207
$ws.Cells.Item(14, 1).NumberFormat = "yyyy-mm-dd"

for ($i = 0; $i -lt $row13.Length; $i++) {
    $col = $i + 2
    $c13 = $ws.Cells.Item(13, $col)
    $c13.Value = $row13[$i]
    $c13.Style = "Normal"

    $c14 = $ws.Cells.Item(14, $col)
    $c14.Value = $row14[$i]
    $c14.Style = "Normal"
}
